$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 7).Value = 26.731658
$ws.Cells.Item(2, 8).Value = 80.194974
$ws.Cells.Item(2, 9).Value = 0.02353393228912
$ws.Cells.Item(2, 10).Value = 0.02353393228912
$ws.Cells.Item(2, 13).Value = 3.339352
$ws.Cells.Item(2, 14).Value = 10.018056
$ws.Cells.Item(2, 15).Value = 0.6054960700393903
$ws.Cells.Item(2, 16).Value = 0.6054960700393903
$ws.Cells.Item(2, 17).Value = 89.266415605616
$ws.Cells.Item(2, 18).Value = 803.3977404505441
$ws.Cells.Item(2, 19).Value = 0.01424970351363527
$ws.Cells.Item(2, 20).Value = 0.01424970351363527
$ws.Cells.Item(3, 7).Value = 26.731658
$ws.Cells.Item(3, 8).Value = 80.194974
$ws.Cells.Item(3, 9).Value = 0.02353393228912
$ws.Cells.Item(3, 10).Value = 0.02353393228912
$ws.Cells.Item(3, 15).Value = 0.2540955070726236
$ws.Cells.Item(3, 16).Value = 0.2540955070726236
$ws.Cells.Item(3, 17).Value = 37.46051586493201
$ws.Cells.Item(3, 18).Value = 337.1446427843881
$ws.Cells.Item(3, 19).Value = 0.005979866458416735
$ws.Cells.Item(3, 20).Value = 0.005979866458416734
$ws.Cells.Item(4, 7).Value = 26.731658
$ws.Cells.Item(4, 8).Value = 80.194974
$ws.Cells.Item(4, 9).Value = 0.02353393228912
$ws.Cells.Item(4, 10).Value = 0.02353393228912
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1338136666666667
$ws.Cells.Item(4, 14).Value = 0.401441
$ws.Cells.Item(4, 15).Value = 0.02426328499787613
$ws.Cells.Item(4, 16).Value = 0.02426328499787612
$ws.Cells.Item(4, 17).Value = 3.577061173059334
$ws.Cells.Item(4, 18).Value = 32.193550557534
$ws.Cells.Item(4, 19).Value = 0.0005710105062516378
$ws.Cells.Item(4, 20).Value = 0.0005710105062516377
$ws.Cells.Item(5, 7).Value = 26.731658
$ws.Cells.Item(5, 8).Value = 80.194974
$ws.Cells.Item(5, 9).Value = 0.02353393228912
$ws.Cells.Item(5, 10).Value = 0.02353393228912
$ws.Cells.Item(5, 13).Value = 0.6405483333333334
$ws.Cells.Item(5, 14).Value = 1.921645
$ws.Cells.Item(5, 15).Value = 0.11614513789011
$ws.Cells.Item(5, 16).Value = 0.11614513789011
$ws.Cells.Item(5, 17).Value = 17.12291897913667
$ws.Cells.Item(5, 18).Value = 154.10627081223
$ws.Cells.Item(5, 19).Value = 0.002733351810816355
$ws.Cells.Item(5, 20).Value = 0.002733351810816354
$ws.Cells.Item(6, 9).Value = 0.9376016087099961
$ws.Cells.Item(6, 10).Value = 0.9376016087099961
$ws.Cells.Item(6, 13).Value = 3.339352
$ws.Cells.Item(6, 14).Value = 10.018056
$ws.Cells.Item(6, 15).Value = 0.6054960700393903
$ws.Cells.Item(6, 16).Value = 0.6054960700393903
$ws.Cells.Item(6, 17).Value = 3556.410966402518
$ws.Cells.Item(6, 18).Value = 32007.69869762266
$ws.Cells.Item(6, 19).Value = 0.5677140893365128
$ws.Cells.Item(6, 20).Value = 0.5677140893365128
$ws.Cells.Item(7, 9).Value = 0.9376016087099961
$ws.Cells.Item(7, 10).Value = 0.9376016087099961
$ws.Cells.Item(7, 15).Value = 0.2540955070726236
$ws.Cells.Item(7, 16).Value = 0.2540955070726236
$ws.Cells.Item(7, 19).Value = 0.2382403561972741
$ws.Cells.Item(7, 20).Value = 0.238240356197274
$ws.Cells.Item(8, 9).Value = 0.9376016087099961
$ws.Cells.Item(8, 10).Value = 0.9376016087099961
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.1338136666666667
$ws.Cells.Item(8, 14).Value = 0.401441
$ws.Cells.Item(8, 15).Value = 0.02426328499787613
$ws.Cells.Item(8, 16).Value = 0.02426328499787612
$ws.Cells.Item(8, 17).Value = 142.5115985340463
$ws.Cells.Item(8, 18).Value = 1282.604386806416
$ws.Cells.Item(8, 19).Value = 0.02274929504659777
$ws.Cells.Item(8, 20).Value = 0.02274929504659777
$ws.Cells.Item(9, 9).Value = 0.9376016087099961
$ws.Cells.Item(9, 10).Value = 0.9376016087099961
$ws.Cells.Item(9, 13).Value = 0.6405483333333334
$ws.Cells.Item(9, 14).Value = 1.921645
$ws.Cells.Item(9, 15).Value = 0.11614513789011
$ws.Cells.Item(9, 16).Value = 0.11614513789011
$ws.Cells.Item(9, 17).Value = 682.1841833917246
$ws.Cells.Item(9, 18).Value = 6139.657650525521
$ws.Cells.Item(9, 19).Value = 0.1088978681296115
$ws.Cells.Item(9, 20).Value = 0.1088978681296115
$ws.Cells.Item(10, 7).Value = 0.1721486666666666
$ws.Cells.Item(10, 8).Value = 0.516446
$ws.Cells.Item(10, 9).Value = 0.0001515556971810586
$ws.Cells.Item(10, 10).Value = 0.0001515556971810586
$ws.Cells.Item(10, 13).Value = 3.339352
$ws.Cells.Item(10, 14).Value = 10.018056
$ws.Cells.Item(10, 15).Value = 0.6054960700393903
$ws.Cells.Item(10, 16).Value = 0.6054960700393903
$ws.Cells.Item(10, 17).Value = 0.5748649943306666
$ws.Cells.Item(10, 18).Value = 5.173784948976
$ws.Cells.Item(10, 19).Value = 0.00009176637903521087
$ws.Cells.Item(10, 20).Value = 0.00009176637903521087
$ws.Cells.Item(11, 7).Value = 0.1721486666666666
$ws.Cells.Item(11, 8).Value = 0.516446
$ws.Cells.Item(11, 9).Value = 0.0001515556971810586
$ws.Cells.Item(11, 10).Value = 0.0001515556971810586
$ws.Cells.Item(11, 15).Value = 0.2540955070726236
$ws.Cells.Item(11, 16).Value = 0.2540955070726236
$ws.Cells.Item(11, 17).Value = 0.241241222628
$ws.Cells.Item(11, 18).Value = 2.171171003652
$ws.Cells.Item(11, 19).Value = 0.00003850962172496607
$ws.Cells.Item(11, 20).Value = 0.00003850962172496607
$ws.Cells.Item(12, 7).Value = 0.1721486666666666
$ws.Cells.Item(12, 8).Value = 0.516446
$ws.Cells.Item(12, 9).Value = 0.0001515556971810586
$ws.Cells.Item(12, 10).Value = 0.0001515556971810586
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.1338136666666667
$ws.Cells.Item(12, 14).Value = 0.401441
$ws.Cells.Item(12, 15).Value = 0.02426328499787613
$ws.Cells.Item(12, 16).Value = 0.02426328499787612
$ws.Cells.Item(12, 17).Value = 0.02303584429844445
$ws.Cells.Item(12, 18).Value = 0.207322598686
$ws.Cells.Item(12, 19).Value = 0.000003677239073755835
$ws.Cells.Item(12, 20).Value = 0.000003677239073755835
$ws.Cells.Item(13, 7).Value = 0.1721486666666666
$ws.Cells.Item(13, 8).Value = 0.516446
$ws.Cells.Item(13, 9).Value = 0.0001515556971810586
$ws.Cells.Item(13, 10).Value = 0.0001515556971810586
$ws.Cells.Item(13, 13).Value = 0.6405483333333334
$ws.Cells.Item(13, 14).Value = 1.921645
$ws.Cells.Item(13, 15).Value = 0.11614513789011
$ws.Cells.Item(13, 16).Value = 0.11614513789011
$ws.Cells.Item(13, 17).Value = 0.1102695415188889
$ws.Cells.Item(13, 18).Value = 0.9924258736699999
$ws.Cells.Item(13, 19).Value = 0.00001760245734712581
$ws.Cells.Item(13, 20).Value = 0.00001760245734712581
$ws.Cells.Item(14, 7).Value = 43.33877
$ws.Cells.Item(14, 8).Value = 130.01631
$ws.Cells.Item(14, 9).Value = 0.03815444888131313
$ws.Cells.Item(14, 10).Value = 0.03815444888131313
$ws.Cells.Item(14, 13).Value = 3.339352
$ws.Cells.Item(14, 14).Value = 10.018056
$ws.Cells.Item(14, 15).Value = 0.6054960700393903
$ws.Cells.Item(14, 16).Value = 0.6054960700393903
$ws.Cells.Item(14, 17).Value = 144.72340827704
$ws.Cells.Item(14, 18).Value = 1302.51067449336
$ws.Cells.Item(14, 19).Value = 0.02310236885215391
$ws.Cells.Item(14, 20).Value = 0.02310236885215391
$ws.Cells.Item(15, 7).Value = 43.33877
$ws.Cells.Item(15, 8).Value = 130.01631
$ws.Cells.Item(15, 9).Value = 0.03815444888131313
$ws.Cells.Item(15, 10).Value = 0.03815444888131313
$ws.Cells.Item(15, 15).Value = 0.2540955070726236
$ws.Cells.Item(15, 16).Value = 0.2540955070726236
$ws.Cells.Item(15, 17).Value = 60.73295869458001
$ws.Cells.Item(15, 18).Value = 546.5966282512201
$ws.Cells.Item(15, 19).Value = 0.009694874035573757
$ws.Cells.Item(15, 20).Value = 0.009694874035573755
$ws.Cells.Item(16, 7).Value = 43.33877
$ws.Cells.Item(16, 8).Value = 130.01631
$ws.Cells.Item(16, 9).Value = 0.03815444888131313
$ws.Cells.Item(16, 10).Value = 0.03815444888131313
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.1338136666666667
$ws.Cells.Item(16, 14).Value = 0.401441
$ws.Cells.Item(16, 15).Value = 0.02426328499787613
$ws.Cells.Item(16, 16).Value = 0.02426328499787612
$ws.Cells.Item(16, 17).Value = 5.799319722523335
$ws.Cells.Item(16, 18).Value = 52.19387750271001
$ws.Cells.Item(16, 19).Value = 0.0009257522671441963
$ws.Cells.Item(16, 20).Value = 0.0009257522671441962
$ws.Cells.Item(17, 7).Value = 43.33877
$ws.Cells.Item(17, 8).Value = 130.01631
$ws.Cells.Item(17, 9).Value = 0.03815444888131313
$ws.Cells.Item(17, 10).Value = 0.03815444888131313
$ws.Cells.Item(17, 13).Value = 0.6405483333333334
$ws.Cells.Item(17, 14).Value = 1.921645
$ws.Cells.Item(17, 15).Value = 0.11614513789011
$ws.Cells.Item(17, 16).Value = 0.11614513789011
$ws.Cells.Item(17, 17).Value = 27.76057689221667
$ws.Cells.Item(17, 18).Value = 249.84519202995
$ws.Cells.Item(17, 19).Value = 0.004431453726441268
$ws.Cells.Item(17, 20).Value = 0.004431453726441268
$ws.Cells.Item(18, 7).Value = 0.4290093333333333
$ws.Cells.Item(18, 8).Value = 1.287028
$ws.Cells.Item(18, 9).Value = 0.0003776898762533613
$ws.Cells.Item(18, 10).Value = 0.0003776898762533613
$ws.Cells.Item(18, 13).Value = 3.339352
$ws.Cells.Item(18, 14).Value = 10.018056
$ws.Cells.Item(18, 15).Value = 0.6054960700393903
$ws.Cells.Item(18, 16).Value = 0.6054960700393903
$ws.Cells.Item(18, 17).Value = 1.432613175285333
$ws.Cells.Item(18, 18).Value = 12.893518577568
$ws.Cells.Item(18, 19).Value = 0.0002286897357650739
$ws.Cells.Item(18, 20).Value = 0.0002286897357650739
$ws.Cells.Item(19, 7).Value = 0.4290093333333333
$ws.Cells.Item(19, 8).Value = 1.287028
$ws.Cells.Item(19, 9).Value = 0.0003776898762533613
$ws.Cells.Item(19, 10).Value = 0.0003776898762533613
$ws.Cells.Item(19, 15).Value = 0.2540955070726236
$ws.Cells.Item(19, 16).Value = 0.2540955070726236
$ws.Cells.Item(19, 17).Value = 0.6011939453040001
$ws.Cells.Item(19, 18).Value = 5.410745507736
$ws.Cells.Item(19, 19).Value = 0.00009596930062279432
$ws.Cells.Item(19, 20).Value = 0.00009596930062279428
$ws.Cells.Item(20, 7).Value = 0.4290093333333333
$ws.Cells.Item(20, 8).Value = 1.287028
$ws.Cells.Item(20, 9).Value = 0.0003776898762533613
$ws.Cells.Item(20, 10).Value = 0.0003776898762533613
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 12).Value = 0.6666666666666666
$ws.Cells.Item(20, 13).Value = 0.1338136666666667
$ws.Cells.Item(20, 14).Value = 0.401441
$ws.Cells.Item(20, 15).Value = 0.02426328499787613
$ws.Cells.Item(20, 16).Value = 0.02426328499787612
$ws.Cells.Item(20, 17).Value = 0.05740731192755556
$ws.Cells.Item(20, 18).Value = 0.5166658073480001
$ws.Cells.Item(20, 19).Value = 0.000009163997108347871
$ws.Cells.Item(20, 20).Value = 0.00000916399710834787
$ws.Cells.Item(21, 7).Value = 0.4290093333333333
$ws.Cells.Item(21, 8).Value = 1.287028
$ws.Cells.Item(21, 9).Value = 0.0003776898762533613
$ws.Cells.Item(21, 10).Value = 0.0003776898762533613
$ws.Cells.Item(21, 13).Value = 0.6405483333333334
$ws.Cells.Item(21, 14).Value = 1.921645
$ws.Cells.Item(21, 15).Value = 0.11614513789011
$ws.Cells.Item(21, 16).Value = 0.11614513789011
$ws.Cells.Item(21, 17).Value = 0.2748012134511111
$ws.Cells.Item(21, 18).Value = 2.47321092106
$ws.Cells.Item(21, 19).Value = 0.00004386684275714525
$ws.Cells.Item(21, 20).Value = 0.00004386684275714523
$ws.Cells.Item(22, 7).Value = 0.2053263333333333
$ws.Cells.Item(22, 8).Value = 0.6159789999999999
$ws.Cells.Item(22, 9).Value = 0.0001807645461362684
$ws.Cells.Item(22, 10).Value = 0.0001807645461362684
$ws.Cells.Item(22, 13).Value = 3.339352
$ws.Cells.Item(22, 14).Value = 10.018056
$ws.Cells.Item(22, 15).Value = 0.6054960700393903
$ws.Cells.Item(22, 16).Value = 0.6054960700393903
$ws.Cells.Item(22, 17).Value = 0.6856569018693333
$ws.Cells.Item(22, 18).Value = 6.170912116824001
$ws.Cells.Item(22, 19).Value = 0.0001094522222879646
$ws.Cells.Item(22, 20).Value = 0.0001094522222879646
$ws.Cells.Item(23, 7).Value = 0.2053263333333333
$ws.Cells.Item(23, 8).Value = 0.6159789999999999
$ws.Cells.Item(23, 9).Value = 0.0001807645461362684
$ws.Cells.Item(23, 10).Value = 0.0001807645461362684
$ws.Cells.Item(23, 15).Value = 0.2540955070726236
$ws.Cells.Item(23, 16).Value = 0.2540955070726236
$ws.Cells.Item(23, 17).Value = 0.287734878522
$ws.Cells.Item(23, 18).Value = 2.589613906698
$ws.Cells.Item(23, 19).Value = 0.00004593145901124779
$ws.Cells.Item(23, 20).Value = 0.00004593145901124778
$ws.Cells.Item(24, 7).Value = 0.2053263333333333
$ws.Cells.Item(24, 8).Value = 0.6159789999999999
$ws.Cells.Item(24, 9).Value = 0.0001807645461362684
$ws.Cells.Item(24, 10).Value = 0.0001807645461362684
$ws.Cells.Item(24, 11).Value = 2
$ws.Cells.Item(24, 12).Value = 0.6666666666666666
$ws.Cells.Item(24, 13).Value = 0.1338136666666667
$ws.Cells.Item(24, 14).Value = 0.401441
$ws.Cells.Item(24, 15).Value = 0.02426328499787613
$ws.Cells.Item(24, 16).Value = 0.02426328499787612
$ws.Cells.Item(24, 17).Value = 0.02747546952655556
$ws.Cells.Item(24, 18).Value = 0.247279225739
$ws.Cells.Item(24, 19).Value = 0.000004385941700416007
$ws.Cells.Item(24, 20).Value = 0.000004385941700416007
$ws.Cells.Item(25, 7).Value = 0.2053263333333333
$ws.Cells.Item(25, 8).Value = 0.6159789999999999
$ws.Cells.Item(25, 9).Value = 0.0001807645461362684
$ws.Cells.Item(25, 10).Value = 0.0001807645461362684
$ws.Cells.Item(25, 13).Value = 0.6405483333333334
$ws.Cells.Item(25, 14).Value = 1.921645
$ws.Cells.Item(25, 15).Value = 0.11614513789011
$ws.Cells.Item(25, 16).Value = 0.11614513789011
$ws.Cells.Item(25, 17).Value = 0.1315214406061111
$ws.Cells.Item(25, 18).Value = 0.9924258736699999
$ws.Cells.Item(25, 19).Value = 0.00002099492313664005
$ws.Cells.Item(25, 20).Value = 0.00002099492313664005
